# Keep most recent assessment in-sync with client changes when uploading
#
# The client's assessment export format changed: several columns were
# dropped (Date First became Homeless, Cumulative Months Homeless in Last
# Three Years, VI-SPDAT Score, Referral Date, Age greater than 65 years of
# age, Minimum Bedroom Size) and replaced with new ones used by the
# updated intake form (Disabled Per HUD Language, Substance Use
# Disability, Mental Health Disability, Cumulative days Homeless, Age
# greater than 60 years of age, Currently first time pregnant 28 weeks or
# less). Rebuild the header row and the sample data row to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Shelter Location"
$ws.Range("B1").Value = "Disabled Per HUD Language"
$ws.Range("C1").Value = "Home-base ID"
$ws.Range("D1").Value = "Substance Use Disability"
$ws.Range("E1").Value = "Mental Health Disability"
$ws.Range("F1").Value = "Occurrences of Homelessness in Last Three Years"
$ws.Range("G1").Value = "Cumulative days Homeless"
$ws.Range("H1").Value = "Family of at least one Adult and one child"
$ws.Range("I1").Value = "Age greater than 60 years of age"
$ws.Range("J1").Value = "Age less than 24 years of age"
$ws.Range("K1").Value = "Permanent Supportive Housing Eligible"
$ws.Range("L1").Value = "Currently first time pregnant 28 weeks or less"
$ws.Range("M1").Value = "Veteran Status"
$ws.Range("N1").Value = "HOPWA Eligible"
$ws.Range("O1").Value = "Prioritized for Health"

# ---- Sample data row (row 2) ----
$ws.Range("A2").Value = "Inside"
$ws.Range("B2").Value = "YES"

$ws.Range("C2").Style = "Normal"
$ws.Range("C2").Value = 5555

$ws.Range("D2").Value = "YES"

$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "YES"

$ws.Range("F2").Style = "Normal"
$ws.Range("F2").Value = 3

$ws.Range("G2").Style = "Normal"
$ws.Range("G2").Value = 200

$ws.Range("H2").Value = "Y"
$ws.Range("I2").Value = "Y"

$ws.Range("J2").Style = "Normal"
$ws.Range("J2").Value = "No"

$ws.Range("K2").Style = "Normal"
$ws.Range("K2").Value = "Y"

$ws.Range("L2").Value = "y"
$ws.Range("M2").Value = "no"

$ws.Range("N2").Style = "Normal"
$ws.Range("N2").Value = "N"

$ws.Range("O2").Value = "Yes"

# ---- Column widths (best-fit widths left behind by Excel after the
#      header text grew) ----
$ws.Columns.Item(2).ColumnWidth = 23.833333333333332
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334
$ws.Columns.Item(4).ColumnWidth = 24.166666666666668
$ws.Columns.Item(5).ColumnWidth = 41.666666666666664
$ws.Columns.Item(6).ColumnWidth = 42.166666666666664
$ws.Columns.Item(7).ColumnWidth = 35.0
$ws.Columns.Item(8).ColumnWidth = 35.5
$ws.Columns.Item(9).ColumnWidth = 28.166666666666668
$ws.Columns.Item(10).ColumnWidth = 32.666666666666664
$ws.Columns.Item(11).ColumnWidth = 33.166666666666664
$ws.Columns.Item(12).ColumnWidth = 39.0
$ws.Columns.Item(13).ColumnWidth = 13.333333333333334
$ws.Columns.Item(14).ColumnWidth = 13.833333333333334
$ws.Columns.Item(15).ColumnWidth = 18.0

# ---- Selection left over from editing the next few (empty) rows ----
$ws.Range("A3:XFD5").Select() | Out-Null
